# Append: 2026-01-14 01:39 JST
# Update the "取得日時" (acquired datetime) column A for all data rows
# (rows 2-18) on the "ランサーズ" sheet from "2026-01-13 18:37:00" to
# "2026-01-14 01:39:58".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-14 01:39:58"

for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
